$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 112
$ws.Range("H112").Value = 17858506
$ws.Range("J112").Value = 21740642
$ws.Range("L112").Value = 65221926
$ws.Range("N112").Value = -65224142
# Row 134
$ws.Range("H134").Value = 37000
$ws.Range("J134").Value = 37000
$ws.Range("L134").Value = 37000
$ws.Range("N134").Value = -47140

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 3304.4546
$ws.Range("I2").Value = 4224.8335
$ws.Range("J2").Value = 2200
$ws.Range("K2").Value = 4224.8335
$ws.Range("L2").Value = 2200
$ws.Range("M2").Value = -4111.8335
$ws.Range("N2").Value = -2426
# Row 32
$ws.Range("H32").Value = 7040.1704
$ws.Range("I32").Value = 6062.0967
$ws.Range("K32").Value = 6062.0967
$ws.Range("M32").Value = -5775.0967
# Row 116
$ws.Range("H116").Value = 3304.4546
$ws.Range("I116").Value = 4224.8335
$ws.Range("J116").Value = 2200
$ws.Range("K116").Value = 4224.8335
$ws.Range("L116").Value = 2200
$ws.Range("M116").Value = -1930.8335
$ws.Range("N116").Value = -6788
# Row 118
$ws.Range("H118").Value = 34796.668
$ws.Range("J118").Value = 34796.668
$ws.Range("L118").Value = 34796.668
$ws.Range("N118").Value = -38110.668

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 3304.4546
$ws.Range("I3").Value = 4224.8335
$ws.Range("J3").Value = 2200
$ws.Range("K3").Value = 4224.8335
$ws.Range("L3").Value = 2200
$ws.Range("M3").Value = -4110.8335
$ws.Range("N3").Value = -2428
# Row 86
$ws.Range("H86").Value = 33335134
$ws.Range("I86").Value = 41668430
$ws.Range("K86").Value = 41668430
$ws.Range("M86").Value = -41667307
# Row 89
$ws.Range("H89").Value = 33335134
$ws.Range("I89").Value = 41668430
$ws.Range("K89").Value = 208342150
$ws.Range("M89").Value = -208336534
# Row 107
$ws.Range("H107").Value = 889.44446
$ws.Range("I107").Value = 920.2222
$ws.Range("J107").Value = 827.8889
$ws.Range("K107").Value = 920.2222
$ws.Range("L107").Value = 827.8889
$ws.Range("M107").Value = 999.7778
$ws.Range("N107").Value = -4667.8889

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 6
$ws.Range("H6").Value = 40000000
$ws.Range("I6").Value = 40000000
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 40000000
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -39999887
$ws.Range("N6").ClearContents()
# Row 7
$ws.Range("H7").Value = 399.55554
$ws.Range("I7").Value = 424.75
$ws.Range("J7").Value = 198
$ws.Range("K7").Value = 424.75
$ws.Range("L7").Value = 198
$ws.Range("M7").Value = -311.75
$ws.Range("N7").Value = -424
# Row 17
$ws.Range("H17").Value = 250041250
$ws.Range("J17").Value = 250041250
$ws.Range("L17").Value = 250041250
$ws.Range("N17").Value = -250041598
# Row 31
$ws.Range("H31").Value = 6532.25
$ws.Range("I31").Value = 1365.1428
$ws.Range("J31").Value = 22033.572
$ws.Range("K31").Value = 1365.1428
$ws.Range("L31").Value = 22033.572
$ws.Range("M31").Value = -1070.1428
$ws.Range("N31").Value = -22623.572
# Row 34
$ws.Range("H34").Value = 6532.25
$ws.Range("I34").Value = 1365.1428
$ws.Range("J34").Value = 22033.572
$ws.Range("K34").Value = 1365.1428
$ws.Range("L34").Value = 22033.572
$ws.Range("M34").Value = -1163.1428
$ws.Range("N34").Value = -22437.572
# Row 41
$ws.Range("H41").Value = 12000
$ws.Range("J41").Value = 17000
$ws.Range("L41").Value = 17000
$ws.Range("N41").Value = -17856
# Row 51
$ws.Range("H51").Value = 18378
$ws.Range("J51").Value = 20472.5
$ws.Range("L51").Value = 20472.5
$ws.Range("N51").Value = -21944.5
# Row 58
$ws.Range("H58").Value = 1179.625
$ws.Range("I58").Value = 994.85
$ws.Range("J58").Value = 2103.5
$ws.Range("K58").Value = 994.85
$ws.Range("L58").Value = 2103.5
$ws.Range("M58").Value = -791.85
$ws.Range("N58").Value = -2509.5
# Row 59
$ws.Range("H59").Value = 25529.5
$ws.Range("J59").Value = 25529.5
$ws.Range("L59").Value = 25529.5
$ws.Range("N59").Value = -27819.5
# Row 60
$ws.Range("H60").Value = 11663.479
$ws.Range("J60").Value = 11663.479
$ws.Range("L60").Value = 11663.479
$ws.Range("N60").Value = -12685.479
# Row 61
$ws.Range("H61").Value = 18378
$ws.Range("J61").Value = 20472.5
$ws.Range("L61").Value = 20472.5
$ws.Range("N61").Value = -21168.5
# Row 68
$ws.Range("H68").Value = 32000
$ws.Range("J68").Value = 32000
$ws.Range("L68").Value = 32000
$ws.Range("N68").Value = -33498
# Row 71
$ws.Range("H71").Value = 32000
$ws.Range("J71").Value = 32000
$ws.Range("L71").Value = 96000
$ws.Range("N71").Value = -103488
# Row 74
$ws.Range("H74").Value = 21231.2
$ws.Range("J74").Value = 21231.2
$ws.Range("L74").Value = 21231.2
$ws.Range("N74").Value = -22979.2
# Row 77
$ws.Range("H77").Value = 21231.2
$ws.Range("J77").Value = 21231.2
$ws.Range("L77").Value = 63693.60000000001
$ws.Range("N77").Value = -72429.60000000001
# Row 132
$ws.Range("H132").Value = 3108.5386
$ws.Range("I132").Value = 2900
$ws.Range("J132").Value = 3351.8333
$ws.Range("K132").Value = 8700
$ws.Range("L132").Value = 10055.4999
$ws.Range("M132").Value = -6170
$ws.Range("N132").Value = -15115.4999
# Row 136
$ws.Range("H136").Value = 1179.625
$ws.Range("I136").Value = 994.85
$ws.Range("J136").Value = 2103.5
$ws.Range("K136").Value = 2984.55
$ws.Range("L136").Value = 6310.5
$ws.Range("M136").Value = -434.5500000000002
$ws.Range("N136").Value = -11410.5

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 17
$ws.Range("H17").Value = 1250
$ws.Range("I17").Value = 1250
$ws.Range("K17").Value = 3750
$ws.Range("M17").Value = -3581
# Row 132
$ws.Range("H132").Value = 3255.1343
$ws.Range("I132").Value = 2560.8
$ws.Range("J132").Value = 3311.1292
$ws.Range("K132").Value = 23047.2
$ws.Range("L132").Value = 29800.1628
$ws.Range("M132").Value = -20517.2
$ws.Range("N132").Value = -34860.1628

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 9
$ws.Range("H9").Value = 47005.332
$ws.Range("J9").Value = 70008
$ws.Range("L9").Value = 70008
$ws.Range("N9").Value = -70348

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 113
$ws.Range("H113").Value = 1622.1818
$ws.Range("I113").Value = 1748.2
$ws.Range("J113").Value = 1517.1666
$ws.Range("K113").Value = 5244.6
$ws.Range("L113").Value = 4551.4998
$ws.Range("M113").Value = -3074.6
$ws.Range("N113").Value = -8891.4998
# Row 123
$ws.Range("H123").Value = 30000
$ws.Range("J123").Value = 30000
$ws.Range("L123").Value = 30000
$ws.Range("N123").Value = -39800
